$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2384290
$ws.Range("I74").Value = 2384290
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2384290
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2383354
$ws.Range("N74").Value = $null

$ws.Range("H75").Value = 29938
$ws.Range("J75").Value = 29938
$ws.Range("L75").Value = 29938
$ws.Range("N75").Value = -31810

$ws.Range("H76").Value = 23812652
$ws.Range("I76").Value = 24393204
$ws.Range("K76").Value = 24393204
$ws.Range("M76").Value = -24392889

$ws.Range("H77").Value = 2384290
$ws.Range("I77").Value = 2384290
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 11921450
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -11916770
$ws.Range("N77").Value = $null

$ws.Range("H78").Value = 29938
$ws.Range("J78").Value = 29938
$ws.Range("L78").Value = 89814
$ws.Range("N78").Value = -99174

$ws.Range("H79").Value = 23812652
$ws.Range("I79").Value = 24393204
$ws.Range("K79").Value = 24393204
$ws.Range("M79").Value = -24392112

$ws.Range("H93").Value = 23750
$ws.Range("J93").Value = 23750
$ws.Range("L93").Value = 23750
$ws.Range("N93").Value = -28742

$ws.Range("H95").Value = 26833.334
$ws.Range("J95").Value = 26833.334
$ws.Range("L95").Value = 26833.334
$ws.Range("N95").Value = -32325.334

$ws.Range("H99").Value = 2780.6365
$ws.Range("I99").Value = 4376.4
$ws.Range("J99").Value = 1450.8334
$ws.Range("K99").Value = 13129.2
$ws.Range("L99").Value = 4352.5002
$ws.Range("M99").Value = -11631.2
$ws.Range("N99").Value = -7348.5002

$ws.Range("H100").Value = 6447.9287
$ws.Range("I100").Value = 4869.5454
$ws.Range("J100").Value = 12235.333
$ws.Range("K100").Value = 4869.5454
$ws.Range("L100").Value = 12235.333
$ws.Range("M100").Value = -4328.5454
$ws.Range("N100").Value = -13317.333

$ws.Range("H101").Value = 2198.125
$ws.Range("I101").Value = 950
$ws.Range("J101").Value = 3446.25
$ws.Range("K101").Value = 2850
$ws.Range("L101").Value = 10338.75
$ws.Range("M101").Value = -1228
$ws.Range("N101").Value = -13582.75

$ws.Range("H103").Value = 312.81818
$ws.Range("I103").Value = 305.8889
$ws.Range("J103").Value = 344
$ws.Range("K103").Value = 917.6667
$ws.Range("L103").Value = 1032
$ws.Range("M103").Value = -331.6667
$ws.Range("N103").Value = -2204

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws.Range("H106").Value = 46907.855
$ws.Range("I106").Value = 46907.855
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 46907.855
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -46276.855
$ws.Range("N106").Value = $null

$ws.Range("H107").Value = 430.83334
$ws.Range("I107").Value = 430.83334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 430.83334
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1489.16666
$ws.Range("N107").Value = $null

$ws.Range("H135").Value = 13515139
$ws.Range("I135").Value = 2034
$ws.Range("J135").Value = 19232222
$ws.Range("K135").Value = 18306
$ws.Range("L135").Value = 173089998
$ws.Range("M135").Value = -15771
$ws.Range("N135").Value = -173095068

$ws.Range("H138").Value = 1678.1698
$ws.Range("I138").Value = 974
$ws.Range("J138").Value = 1982.6757
$ws.Range("K138").Value = 2922
$ws.Range("L138").Value = 5948.0271
$ws.Range("M138").Value = 2218
$ws.Range("N138").Value = -16228.0271

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3886.7144
$ws.Range("I61").Value = 2117.8333
$ws.Range("K61").Value = 2117.8333
$ws.Range("M61").Value = -1905.8333

$ws.Range("H88").Value = 15075.177
$ws.Range("I88").Value = 5496
$ws.Range("J88").Value = 19066.5
$ws.Range("K88").Value = 5496
$ws.Range("L88").Value = 19066.5
$ws.Range("M88").Value = -5090
$ws.Range("N88").Value = -19878.5

$ws.Range("H91").Value = 15075.177
$ws.Range("I91").Value = 5496
$ws.Range("J91").Value = 19066.5
$ws.Range("K91").Value = 5496
$ws.Range("L91").Value = 19066.5
$ws.Range("M91").Value = -4092
$ws.Range("N91").Value = -21874.5

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

$ws.Range("H132").Value = 14395621
$ws.Range("I132").Value = 23954692
$ws.Range("J132").Value = 1012922
$ws.Range("K132").Value = 71864076
$ws.Range("L132").Value = 3038766
$ws.Range("M132").Value = -71861546
$ws.Range("N132").Value = -3043826

$ws.Range("H136").Value = 3886.7144
$ws.Range("I136").Value = 2117.8333
$ws.Range("K136").Value = 6353.499899999999
$ws.Range("M136").Value = -3803.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 22099.5
$ws.Range("J81").Value = 22099.5
$ws.Range("L81").Value = 22099.5
$ws.Range("N81").Value = -24221.5

$ws.Range("H84").Value = 22099.5
$ws.Range("J84").Value = 22099.5
$ws.Range("L84").Value = 66298.5
$ws.Range("N84").Value = -76906.5

$ws.Range("H86").Value = 501764.28
$ws.Range("I86").Value = 1641.5834
$ws.Range("K86").Value = 1641.5834
$ws.Range("M86").Value = -518.5834

$ws.Range("H89").Value = 501764.28
$ws.Range("I89").Value = 1641.5834
$ws.Range("K89").Value = 8207.916999999999
$ws.Range("M89").Value = -2591.916999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 100003280
$ws.Range("I62").Value = 100003280
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 100003280
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -100002656
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 100003280
$ws.Range("I65").Value = 100003280
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 500016400
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -500013280
$ws.Range("N65").Value = $null

$ws.Range("H132").Value = 2260.5
$ws.Range("I132").Value = 1261.5
$ws.Range("J132").Value = 3059.7
$ws.Range("K132").Value = 3784.5
$ws.Range("L132").Value = 9179.099999999999
$ws.Range("M132").Value = -1254.5
$ws.Range("N132").Value = -14239.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 360.66666
$ws.Range("J122").Value = 466
$ws.Range("L122").Value = 4194
$ws.Range("N122").Value = -9094

$ws.Range("H126").Value = 3536.6667
$ws.Range("I126").Value = 2715
$ws.Range("K126").Value = 8145
$ws.Range("M126").Value = -3205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 48450.816
$ws.Range("I132").Value = 2302
$ws.Range("J132").Value = 65756.625
$ws.Range("K132").Value = 6906
$ws.Range("L132").Value = 197269.875
$ws.Range("M132").Value = -4376
$ws.Range("N132").Value = -202329.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524

$ws.Range("H132").Value = 47030.49
$ws.Range("I132").Value = 68511.47
$ws.Range("J132").Value = 4068.5334
$ws.Range("K132").Value = 205534.41
$ws.Range("L132").Value = 12205.6002
$ws.Range("M132").Value = -203004.41
$ws.Range("N132").Value = -17265.6002
